# Update "想去人数" (want-to-go count) figures for four events on both
# the "展览" sheet and the "全部类型" sheet (which mirrors the same data).
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 1189
    $ws.Range("F3").Value = 602
    $ws.Range("F10").Value = 5539
    $ws.Range("F11").Value = 4922
}
